{"js": "// Update the date heading and every \"a op b=\" cell in the 20x5 practice\n// table to the new values from the target revision. Each table cell\n// holds exactly one paragraph/run, so we replace the range text in place\n// (Word.InsertLocation.replace) rather than inserting a brand new run -\n// this preserves the existing run/paragraph formatting (fonts, size,\n// alignment) instead of losing it.\n\nconst NEW_VALUES = {\"date\":\"2023-05-04 Thursday\",\"grid\":[[\"22-16=\",\"19+33=\",\"27+58=\",\"63+6=\",\"70+0=\"],[\"99-33=\",\"12+69=\",\"55+9=\",\"54-32=\",\"55+22=\"],[\"99-76=\",\"90-73=\",\"9+53=\",\"28+42=\",\"19-3=\"],[\"51-28=\",\"75-72=\",\"6+83=\",\"39-30=\",\"95-57=\"],[\"42+19=\",\"93-48=\",\"11+0=\",\"9-1=\",\"67+2=\"],[\"14+78=\",\"8+85=\",\"26+23=\",\"51-47=\",\"80-2=\"],[\"67+14=\",\"21+61=\",\"0+96=\",\"94-94=\",\"38+21=\"],[\"59+4=\",\"90-65=\",\"70-31=\",\"10-2=\",\"50-19=\"],[\"95-35=\",\"83-39=\",\"52+14=\",\"67+27=\",\"87-62=\"],[\"12+80=\",\"57-33=\",\"74+14=\",\"35+45=\",\"38+8=\"],[\"6+55=\",\"16+45=\",\"64-8=\",\"68-42=\",\"94-10=\"],[\"7+30=\",\"84+13=\",\"41-9=\",\"38-5=\",\"76-46=\"],[\"51-37=\",\"90+1=\",\"82-47=\",\"30-20=\",\"75+14=\"],[\"7+56=\",\"92-91=\",\"75-64=\",\"20+45=\",\"0+5=\"],[\"7+25=\",\"55-20=\",\"28+40=\",\"4+20=\",\"83-31=\"],[\"59+5=\",\"69+19=\",\"1+30=\",\"81-32=\",\"59+18=\"],[\"88-38=\",\"44+27=\",\"25+52=\",\"49+39=\",\"54-33=\"],[\"11+65=\",\"86+7=\",\"29+8=\",\"87-77=\",\"97-34=\"],[\"90-25=\",\"28+51=\",\"98-5=\",\"52+35=\",\"48-44=\"],[\"7+10=\",\"86-72=\",\"54+13=\",\"74-69=\",\"70-56=\"]]};\n\n// 1) Update the date paragraph above the table (first paragraph in the body).\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nparagraphs.items[0].getRange().insertText(NEW_VALUES.date, Word.InsertLocation.replace);\n\n// 2) Update every cell of the (only) table, in row-major order, matching\n// the order the original/target cell text appears in the document.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\nfor (let r = 0; r < NEW_VALUES.grid.length; r++) {\n  const row = NEW_VALUES.grid[r];\n  for (let c = 0; c < row.length; c++) {\n    const cell = table.getCell(r, c);\n    cell.body.getRange().insertText(row[c], Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) Update the date heading (first paragraph in the document body).\n$d.Paragraphs(1).Range.Text = \"2023-05-04 Thursday\"\n\n# 2) Update every cell of the (only) table, row by row, column by column,\n#    matching the order the original/target cell text appears in the document.\n#    Using Cell().Range.Text = \"...\" keeps the existing run/paragraph\n#    formatting (fonts, size, alignment) instead of inserting a new run.\n$newValues = @(\n    @(\"22-16=\", \"19+33=\", \"27+58=\", \"63+6=\", \"70+0=\"),\n    @(\"99-33=\", \"12+69=\", \"55+9=\", \"54-32=\", \"55+22=\"),\n    @(\"99-76=\", \"90-73=\", \"9+53=\", \"28+42=\", \"19-3=\"),\n    @(\"51-28=\", \"75-72=\", \"6+83=\", \"39-30=\", \"95-57=\"),\n    @(\"42+19=\", \"93-48=\", \"11+0=\", \"9-1=\", \"67+2=\"),\n    @(\"14+78=\", \"8+85=\", \"26+23=\", \"51-47=\", \"80-2=\"),\n    @(\"67+14=\", \"21+61=\", \"0+96=\", \"94-94=\", \"38+21=\"),\n    @(\"59+4=\", \"90-65=\", \"70-31=\", \"10-2=\", \"50-19=\"),\n    @(\"95-35=\", \"83-39=\", \"52+14=\", \"67+27=\", \"87-62=\"),\n    @(\"12+80=\", \"57-33=\", \"74+14=\", \"35+45=\", \"38+8=\"),\n    @(\"6+55=\", \"16+45=\", \"64-8=\", \"68-42=\", \"94-10=\"),\n    @(\"7+30=\", \"84+13=\", \"41-9=\", \"38-5=\", \"76-46=\"),\n    @(\"51-37=\", \"90+1=\", \"82-47=\", \"30-20=\", \"75+14=\"),\n    @(\"7+56=\", \"92-91=\", \"75-64=\", \"20+45=\", \"0+5=\"),\n    @(\"7+25=\", \"55-20=\", \"28+40=\", \"4+20=\", \"83-31=\"),\n    @(\"59+5=\", \"69+19=\", \"1+30=\", \"81-32=\", \"59+18=\"),\n    @(\"88-38=\", \"44+27=\", \"25+52=\", \"49+39=\", \"54-33=\"),\n    @(\"11+65=\", \"86+7=\", \"29+8=\", \"87-77=\", \"97-34=\"),\n    @(\"90-25=\", \"28+51=\", \"98-5=\", \"52+35=\", \"48-44=\"),\n    @(\"7+10=\", \"86-72=\", \"54+13=\", \"74-69=\", \"70-56=\"),\n)\n\n$table = $d.Tables(1)\nfor ($r = 0; $r -lt $newValues.Length; $r++) {\n    $row = $newValues[$r]\n    for ($c = 0; $c -lt $row.Length; $c++) {\n        $table.Cell($r + 1, $c + 1).Range.Text = $row[$c]\n    }\n}\n"}
